$wb = $excel.ActiveWorkbook

$ws_sheet1 = $wb.Worksheets.Item("展览")
$ws_sheet1.Range("F4").Value = 256
$ws_sheet1.Range("F5").Value = 2847
$ws_sheet1.Range("F6").Value = 61
$ws_sheet1.Range("F8").Value = 2193
$ws_sheet1.Range("F9").Value = 305
$ws_sheet1.Range("F13").Value = 2529
$ws_sheet1.Range("F15").Value = 1306
$ws_sheet1.Range("F16").Value = 4623
$ws_sheet1.Range("F18").Value = 4943
$ws_sheet1.Range("F19").Value = 1503
$ws_sheet1.Range("F20").Value = 2832
$ws_sheet1.Range("F21").Value = 3226
$ws_sheet1.Range("F23").Value = 1524
$ws_sheet1.Range("F24").Value = 244
$ws_sheet1.Range("F25").Value = 831
$ws_sheet1.Range("F26").Value = 95
$ws_sheet1.Range("F27").Value = 275
$ws_sheet1.Range("F28").Value = 934
$ws_sheet1.Range("F29").Value = 1720
$ws_sheet1.Range("F30").Value = 112
$ws_sheet1.Range("F31").Value = 264
$ws_sheet1.Range("F32").Value = 652
$ws_sheet1.Range("F33").Value = 153
$ws_sheet1.Range("F34").Value = 316
$ws_sheet1.Range("F35").Value = 385

$ws_sheet2 = $wb.Worksheets.Item("演出")
$ws_sheet2.Range("F3").Value = 96
$ws_sheet2.Range("F8").Value = 86

$ws_sheet4 = $wb.Worksheets.Item("全部类型")
$ws_sheet4.Range("F3").Value = 96
$ws_sheet4.Range("F10").Value = 256
$ws_sheet4.Range("F11").Value = 2847
$ws_sheet4.Range("F12").Value = 61
$ws_sheet4.Range("F13").Value = 2193
$ws_sheet4.Range("F14").Value = 305
$ws_sheet4.Range("F15").Value = 86
$ws_sheet4.Range("F21").Value = 2529
$ws_sheet4.Range("F22").Value = 1306
$ws_sheet4.Range("F26").Value = 4623
$ws_sheet4.Range("F28").Value = 4943
$ws_sheet4.Range("F29").Value = 1503
$ws_sheet4.Range("F30").Value = 2832
$ws_sheet4.Range("F31").Value = 3226
$ws_sheet4.Range("F35").Value = 1524
$ws_sheet4.Range("F37").Value = 244
$ws_sheet4.Range("F38").Value = 831
$ws_sheet4.Range("F39").Value = 95
$ws_sheet4.Range("F40").Value = 275
$ws_sheet4.Range("F41").Value = 934
$ws_sheet4.Range("F43").Value = 1720
$ws_sheet4.Range("F44").Value = 112
$ws_sheet4.Range("F45").Value = 264
$ws_sheet4.Range("F46").Value = 652
$ws_sheet4.Range("F47").Value = 153
$ws_sheet4.Range("F48").Value = 316
$ws_sheet4.Range("F49").Value = 385
